$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1073.8
$ws.Range("J17").Value = 1073.8
$ws.Range("L17").Value = 3221.4
$ws.Range("N17").Value = -3557.4
$ws.Range("H40").Value = 1786.6
$ws.Range("J40").Value = 1800
$ws.Range("L40").Value = 1800
$ws.Range("N40").Value = -2150
$ws.Range("H125").Value = 3745.4285
$ws.Range("I125").Value = 5170.5
$ws.Range("J125").Value = 1845.3334
$ws.Range("K125").Value = 46534.5
$ws.Range("L125").Value = 16608.0006
$ws.Range("M125").Value = -44074.5
$ws.Range("N125").Value = -21528.0006
$ws.Range("H137").Value = 1370.1
$ws.Range("I137").Value = 989.1579
$ws.Range("K137").Value = 2967.4737
$ws.Range("M137").Value = -417.4737

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1276
$ws.Range("I74").Value = 677.3043
$ws.Range("J74").Value = 2997.25
$ws.Range("K74").Value = 677.3043
$ws.Range("L74").Value = 2997.25
$ws.Range("M74").Value = 196.6957
$ws.Range("N74").Value = -4745.25
$ws.Range("H77").Value = 1276
$ws.Range("I77").Value = 677.3043
$ws.Range("J77").Value = 2997.25
$ws.Range("K77").Value = 3386.5215
$ws.Range("L77").Value = 14986.25
$ws.Range("M77").Value = 981.4785000000002
$ws.Range("N77").Value = -23722.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 252474460
$ws.Range("I105").Value = 336631300
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 336631300
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -336629553
$ws.Range("N105").Value = -7494

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 100214.57
$ws.Range("I22").Value = 260
$ws.Range("K22").Value = 260
$ws.Range("M22").Value = 90
$ws.Range("H31").Value = 1352.2258
$ws.Range("I31").Value = 899.4706
$ws.Range("J31").Value = 1902
$ws.Range("K31").Value = 899.4706
$ws.Range("L31").Value = 1902
$ws.Range("M31").Value = -604.4706
$ws.Range("N31").Value = -2492
$ws.Range("H34").Value = 1352.2258
$ws.Range("I34").Value = 899.4706
$ws.Range("J34").Value = 1902
$ws.Range("K34").Value = 899.4706
$ws.Range("L34").Value = 1902
$ws.Range("M34").Value = -697.4706
$ws.Range("N34").Value = -2306
$ws.Range("H62").Value = 9093241
$ws.Range("I62").Value = 2442.8572
$ws.Range("J62").Value = 200000000
$ws.Range("K62").Value = 2442.8572
$ws.Range("L62").Value = 200000000
$ws.Range("M62").Value = -1818.8572
$ws.Range("N62").Value = -200001248
$ws.Range("H65").Value = 9093241
$ws.Range("I65").Value = 2442.8572
$ws.Range("J65").Value = 200000000
$ws.Range("K65").Value = 12214.286
$ws.Range("L65").Value = 1000000000
$ws.Range("M65").Value = -9094.286
$ws.Range("N65").Value = -1000006240

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 7051.467
$ws.Range("I107").Value = 285.83334
$ws.Range("K107").Value = 857.5000200000001
$ws.Range("M107").Value = 1062.49998
$ws.Range("H127").Value = 5500
$ws.Range("J127").Value = 5500
$ws.Range("L127").Value = 16500
$ws.Range("N127").Value = -26420
$ws.Range("H131").Value = 10001998
$ws.Range("I131").Value = 166666930
$ws.Range("J131").Value = 2109.4893
$ws.Range("K131").Value = 500000790
$ws.Range("L131").Value = 6328.467900000001
$ws.Range("M131").Value = -499995750
$ws.Range("N131").Value = -16408.4679

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 56265250
$ws.Range("I70").Value = 41686332
$ws.Range("J70").Value = 100001990
$ws.Range("K70").Value = 41686332
$ws.Range("L70").Value = 100001990
$ws.Range("M70").Value = -41686062
$ws.Range("N70").Value = -100002530
$ws.Range("H73").Value = 56265250
$ws.Range("I73").Value = 41686332
$ws.Range("J73").Value = 100001990
$ws.Range("K73").Value = 41686332
$ws.Range("L73").Value = 100001990
$ws.Range("M73").Value = -41685396
$ws.Range("N73").Value = -100003862
$ws.Range("H80").Value = 4457.143
$ws.Range("I80").Value = 2533.3333
$ws.Range("J80").Value = 5900
$ws.Range("K80").Value = 2533.3333
$ws.Range("L80").Value = 5900
$ws.Range("M80").Value = -1535.3333
$ws.Range("N80").Value = -7896
$ws.Range("H83").Value = 4457.143
$ws.Range("I83").Value = 2533.3333
$ws.Range("J83").Value = 5900
$ws.Range("K83").Value = 12666.6665
$ws.Range("L83").Value = 29500
$ws.Range("M83").Value = -7674.666499999999
$ws.Range("N83").Value = -39484
$ws.Range("H97").Value = 805.1667
$ws.Range("I97").Value = 744
$ws.Range("J97").Value = 1111
$ws.Range("K97").Value = 744
$ws.Range("L97").Value = 1111
$ws.Range("M97").Value = -248
$ws.Range("N97").Value = -2103
$ws.Range("H132").Value = 2089.5173
$ws.Range("I132").Value = 1328.2222
$ws.Range("J132").Value = 3335.2727
$ws.Range("K132").Value = 3984.6666
$ws.Range("L132").Value = 10005.8181
$ws.Range("M132").Value = -1454.6666
$ws.Range("N132").Value = -15065.8181

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 614.7778
$ws.Range("J16").Value = 555
$ws.Range("L16").Value = 555
$ws.Range("N16").Value = -895
$ws.Range("H46").Value = 1666.3334
$ws.Range("I46").Value = 999.5
$ws.Range("K46").Value = 999.5
$ws.Range("M46").Value = -811.5
$ws.Range("H68").Value = 2188.625
$ws.Range("I68").Value = 1801.8182
$ws.Range("J68").Value = 3039.6
$ws.Range("K68").Value = 1801.8182
$ws.Range("L68").Value = 3039.6
$ws.Range("M68").Value = -1052.8182
$ws.Range("N68").Value = -4537.6
$ws.Range("H71").Value = 2188.625
$ws.Range("I71").Value = 1801.8182
$ws.Range("J71").Value = 3039.6
$ws.Range("K71").Value = 9009.091
$ws.Range("L71").Value = 15198
$ws.Range("M71").Value = -5265.091
$ws.Range("N71").Value = -22686
$ws.Range("H93").Value = 880.1
$ws.Range("I93").Value = 799.6667
$ws.Range("K93").Value = 799.6667
$ws.Range("M93").Value = 448.3333
$ws.Range("H132").Value = 20773.424
$ws.Range("I132").Value = 965.5
$ws.Range("K132").Value = 2896.5
$ws.Range("M132").Value = -366.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1606.56
$ws.Range("I132").Value = 1265.2858
$ws.Range("K132").Value = 3795.8574
$ws.Range("M132").Value = -1265.8574
